$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10420603
$ws.Range("I51").Value = 3995.6667
$ws.Range("K51").Value = 3995.6667
$ws.Range("M51").Value = -3511.6667

$ws.Range("H62").Value = 7663.5
$ws.Range("I62").Value = 5498
$ws.Range("J62").Value = 11994.5
$ws.Range("K62").Value = 5498
$ws.Range("L62").Value = 11994.5
$ws.Range("M62").Value = -4874
$ws.Range("N62").Value = -13242.5

$ws.Range("H65").Value = 7663.5
$ws.Range("I65").Value = 5498
$ws.Range("J65").Value = 11994.5
$ws.Range("K65").Value = 27490
$ws.Range("L65").Value = 59972.5
$ws.Range("M65").Value = -24370
$ws.Range("N65").Value = -66212.5

$ws.Range("H98").Value = 2557.7856
$ws.Range("I98").Value = 1317.5834
$ws.Range("J98").Value = 9999
$ws.Range("K98").Value = 1317.5834
$ws.Range("L98").Value = 9999
$ws.Range("M98").Value = 180.4166
$ws.Range("N98").Value = -12995

$ws.Range("H99").Value = 333333440
$ws.Range("I99").Value = 149.5
$ws.Range("J99").Value = 1000000000
$ws.Range("K99").Value = 448.5
$ws.Range("L99").Value = 3000000000
$ws.Range("M99").Value = 1049.5
$ws.Range("N99").Value = -3000002996

$ws.Range("H106").Value = 12849.3
$ws.Range("I106").Value = 2812.25
$ws.Range("K106").Value = 2812.25
$ws.Range("M106").Value = -2181.25

$ws.Range("H113").Value = 49060.363
$ws.Range("J113").Value = 4207.643
$ws.Range("L113").Value = 4207.643
$ws.Range("N113").Value = -10715.643

$ws.Range("H118").Value = 1000000
$ws.Range("I118").Value = 1000000
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 3000000
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -2998343
$ws.Range("N118").ClearContents()

$ws.Range("H122").Value = 2557.7856
$ws.Range("I122").Value = 1317.5834
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 3952.7502
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -1502.7502
$ws.Range("N122").Value = -34897

$ws.Range("H132").Value = 1580.1428
$ws.Range("I132").Value = 1587.0975
$ws.Range("K132").Value = 4761.2925
$ws.Range("M132").Value = -2231.2925

$ws.Range("H137").Value = 1131
$ws.Range("I137").Value = 1053.5555
$ws.Range("K137").Value = 3160.6665
$ws.Range("M137").Value = -610.6664999999998

$ws.Range("H138").Value = 2259.682
$ws.Range("J138").Value = 2556.25
$ws.Range("L138").Value = 7668.75
$ws.Range("N138").Value = -17948.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1918
$ws.Range("I122").Value = 1776.7142
$ws.Range("K122").Value = 5330.142599999999
$ws.Range("M122").Value = -2880.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H47").Value = 126970
$ws.Range("J47").Value = 126970
$ws.Range("L47").Value = 126970
$ws.Range("N47").Value = -128010

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 156.5
$ws.Range("I22").Value = 156.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 156.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 193.5
$ws.Range("N22").ClearContents()

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H132").Value = 4305.5884
$ws.Range("I132").Value = 4265.5
$ws.Range("K132").Value = 12796.5
$ws.Range("M132").Value = -10266.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 629.53845
$ws.Range("I34").Value = 179.125
$ws.Range("J34").Value = 1350.2
$ws.Range("K34").Value = 537.375
$ws.Range("L34").Value = 4050.6
$ws.Range("M34").Value = -453.375
$ws.Range("N34").Value = -4218.6

$ws.Range("H50").Value = 290.7143
$ws.Range("J50").Value = 101
$ws.Range("L50").Value = 303
$ws.Range("N50").Value = -1265

$ws.Range("H53").Value = 290.7143
$ws.Range("J53").Value = 101
$ws.Range("L53").Value = 303
$ws.Range("N53").Value = -1265

$ws.Range("H86").Value = 638.1739
$ws.Range("I86").Value = 450
$ws.Range("K86").Value = 1350
$ws.Range("M86").Value = -164

$ws.Range("H89").Value = 638.1739
$ws.Range("I89").Value = 450
$ws.Range("K89").Value = 4050
$ws.Range("M89").Value = 1878

$ws.Range("H114").Value = 50001676
$ws.Range("I114").Value = 66667230
$ws.Range("K114").Value = 200001690
$ws.Range("M114").Value = -199998436

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 37076.77
$ws.Range("J46").Value = 34999.6
$ws.Range("L46").Value = 34999.6
$ws.Range("N46").Value = -35311.6

$ws.Range("H70").Value = 8411.462
$ws.Range("I70").Value = 8304.637000000001
$ws.Range("K70").Value = 8304.637000000001
$ws.Range("M70").Value = -8034.637000000001

$ws.Range("H73").Value = 8411.462
$ws.Range("I73").Value = 8304.637000000001
$ws.Range("K73").Value = 8304.637000000001
$ws.Range("M73").Value = -7368.637000000001

$ws.Range("H122").Value = 1076.7273
$ws.Range("I122").Value = 1076.7273
$ws.Range("K122").Value = 3230.1819
$ws.Range("M122").Value = -780.1819

$ws.Range("H135").Value = 83309.164
$ws.Range("J135").Value = 84086.13
$ws.Range("L135").Value = 84086.13
$ws.Range("N135").Value = -94226.13

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7840.6665
$ws.Range("I40").Value = 7934.074
$ws.Range("K40").Value = 7934.074
$ws.Range("M40").Value = -7798.074

$ws.Range("H46").Value = 24395.79
$ws.Range("I46").Value = 71734
$ws.Range("J46").Value = 2547.3845
$ws.Range("K46").Value = 71734
$ws.Range("L46").Value = 2547.3845
$ws.Range("M46").Value = -71546
$ws.Range("N46").Value = -2923.3845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 18250
$ws.Range("I2").Value = 18250
$ws.Range("K2").Value = 18250
$ws.Range("M2").Value = -18138

$ws.Range("H14").Value = 1216900
$ws.Range("I14").Value = 3001250
$ws.Range("J14").Value = 27333.334
$ws.Range("K14").Value = 3001250
$ws.Range("L14").Value = 27333.334
$ws.Range("M14").Value = -3001082
$ws.Range("N14").Value = -27669.334

$ws.Range("H113").Value = 493.42856
$ws.Range("I113").Value = 163.75
$ws.Range("K113").Value = 491.25
$ws.Range("M113").Value = 1678.75

$ws.Range("H132").Value = 2145.0715
$ws.Range("I132").Value = 1103.8572
$ws.Range("K132").Value = 3311.5716
$ws.Range("M132").Value = -781.5715999999998

$ws.Range("H136").Value = 3448.6
$ws.Range("I136").Value = 2776.2222
$ws.Range("K136").Value = 8328.6666
$ws.Range("M136").Value = -5778.6666

Write-Host "Applied Leviathan_Profits updates"
